$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("B21").Value = "git remote add origin https://github.com/Angela1094/CursoGit.git + tecla enter"
$ws.Range("C21").Value = "Subir el repositorio a Github"

# Row 22
$ws.Range("B22").Value = "`$ git push origin master "
$ws.Range("C22").Value = "Si no sale con el comando anterior el usuario y contraseña para subir el repositorio. Sirve para actualizar los repositorios"

# Row 23
$ws.Range("A23").Value = "Video 4"
$ws.Range("B23").Value = "git pull"
$ws.Range("C23").Value = "Detecta cambios que hemos hecho en GitHub"

# Wrap text for the two "command" cells in the new block (B21, B22)
$ws.Range("B21:B22").WrapText = $true

# Row heights to match the wrapped content (row 23 keeps the natural/default height)
$ws.Rows.Item(21).RowHeight = 35.05
$ws.Rows.Item(22).RowHeight = 13.8

# Scroll / selection to match the new active cell
$ws.Range("C23").Select()
$excel.ActiveWindow.ScrollRow = 13
